# Regenerate save_data column "K" (column G) values for szapucki_thomas.xlsx
# The commit message indicates the author switched from using the raw
# "Strike#" figure to the computed "K" value, and rewrote the G column
# (rows 2-12) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 3
    9  = 2
    10 = 4
    11 = 1
    12 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
